$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee Info")
$ws.Rows.Item(59).Delete()
